$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Estado" column for rows 13-16: "Pendiente" -> "Completado " ---
$ws.Range("E13").Value = "Completado "
$ws.Range("E14").Value = "Completado "
$ws.Range("E15").Value = "Completado "
$ws.Range("E16").Value = "Completado "

# --- Build new rows 17 and 18 (previously blank placeholder rows) ---
# Copy formatting from existing data rows so borders/alignment/number
# formats match the rest of the table, then overwrite with the new values.
$ws.Range("B13:E13").Copy()
$ws.Range("B17:E17").PasteSpecial(-4122)
$ws.Range("B16:E16").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B17").Value = 11
$ws.Range("C17").Value = "Corrección del ingreso de datos para el registro de exámenes"
$ws.Range("D17").Value = "Backend"
$ws.Range("E17").Value = "Completado "

$ws.Range("B18").Value = 12
$ws.Range("C18").Value = "Asignación de alternativas a cada pregunta creada por examen según su orden y cantidad"
$ws.Range("D18").Value = "Backend "
$ws.Range("E18").Value = "Completado "

# Match the row heights used by the rest of the wrapped-text table rows.
$ws.Range("B17:E17").RowHeight = 28.8
$ws.Range("B18:E18").RowHeight = 28.8

# --- Update view: scroll down a bit and select the newly added rows ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B17:E18").Select()
